$wb = $excel.ActiveWorkbook

# Rename sheets to uppercase (matching target naming in the diff)
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the now-unused "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Keep "PAINEIS DARQ" as the active/selected sheet (it was the active tab originally)
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
